# Robotron 2084 "Enforcer" sprite map — recolor the "G" (capital) cells that
# were using shared-string "G" under the grey style (s=8) to a new lower-case
# "g" label, matching the author's palette tweak. The other "G" cells that use
# the dedicated style s=11 (the actual Green swatches) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @(
    "F2",
    "E3","F3","G3",
    "D4","H4",
    "C5","I5",
    "E6","F6","G6","AB6","AU6",
    "D7","E7","G7","H7","AA7","AC7","AT7","AV7",
    "AU8","BK8",
    "F10",
    "AB19","AU19","AV19",
    "AA20","AC20","AT20","AW20",
    "F21","G21","Z21","AD21","AS21","AX21",
    "E22","H22","AB22","AU22","AV22",
    "F23","G23","AA23","AC23","AT23","AU23","AV23","AW23",
    "AB26","AU26","AV26"
)

foreach ($addr in $cells) {
    $ws.Range($addr).Value = "g"
}

# Reflect the author's last on-screen state: zoomed out a bit and the
# selection left on AP2 (the mirrored index column for the second half of
# the sprite strip).
$ws.Range("AP2").Select()
$excel.ActiveWindow.Zoom = 115
